$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "72242805"
$ws.Range("D16").Value = "RONALD ALBERTO MARQUEZ BARRAZA"
$ws.Range("E16").Value = "1912"
$ws.Range("F16").Value = 28708
$ws.Range("G16").Value = 877803

$ws.Range("C17").Value = "1143332689"
$ws.Range("D17").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E17").Value = "2211"
$ws.Range("F17").Value = 17333
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "1143332689"
$ws.Range("D18").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E18").Value = "2212"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

$ws.Range("C19").Value = "1143332689"
$ws.Range("D19").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

$ws.Range("C20").Value = "1143332689"
$ws.Range("D20").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

$ws.Range("C21").Value = "1143332689"
$ws.Range("D21").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E21").Value = "2303"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

$ws.Range("C22").Value = "1143332689"
$ws.Range("D22").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E22").Value = "2304"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1000000

$ws.Range("C23").Value = "1143332689"
$ws.Range("D23").Value = "GERARDO RAFAEL ACUÑA GONZALEZ"
$ws.Range("E23").Value = "2305"
$ws.Range("F23").Value = 14667
$ws.Range("G23").Value = 1000000
